$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Selection changes on pre-existing sheets
# ------------------------------------------------------------------
$wsOrganMaleMass = $wb.Worksheets.Item("Organ Mass - Male")
$wsOrganMaleMass.Range("C25").Select()

$wsOrganComposition = $wb.Worksheets.Item("Organ Composition")
$wsOrganComposition.Range("A3:B16").Select()

# ------------------------------------------------------------------
# 2) Add the new "Insulin Receptors" sheet as the last tab
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Insulin Receptors"

# ------------------------------------------------------------------
# 3) Fill in the sheet content. Cells are written in the same order
#    the original author entered them so new shared-string entries
#    land at the same indices as the target workbook.
# ------------------------------------------------------------------

$ws.Range("A1").Value = "Insulin Receptors (Units are mU, mU/mL and mU/Min)"

$ws.Range("A3").Value = "Body Mass"
$ws.Range("B3").Formula = "=BodyMassMale"

$ws.Range("A5").Value = "Organ"
$ws.Range("B5").Value = "Mass (G)"

$ws.Range("A6").Value = "Bone"
$ws.Range("B6").Formula = "=BoneMassMale"

$ws.Range("A7").Value = "Brain"
$ws.Range("B7").Formula = "=BrainMassMale"

$ws.Range("A8").Value = "Fat"
$ws.Range("B8").Formula = "=FatMassMale"

$ws.Range("A9").Value = "GITract"
$ws.Range("B9").Formula = "=GIMassMale"

$ws.Range("A10").Value = "Kidney"
$ws.Range("B10").Formula = "=KidneyMassMale"

$ws.Range("A11").Value = "Left Heart"
$ws.Range("B11").Formula = "=LHeartMassMale"

$ws.Range("A12").Value = "Liver"
$ws.Range("B12").Formula = "=LiverMassMale"

$ws.Range("A13").Value = "Other Tissue"
$ws.Range("B13").Formula = "=OtherMassMale"

$ws.Range("A14").Value = "Respiratory Muscle"
$ws.Range("B14").Formula = "=RMuscleMassMale"

$ws.Range("A15").Value = "Right Heart"
$ws.Range("B15").Formula = "=RHeartMassMale"

$ws.Range("A16").Value = "Skeletal Muscle"
$ws.Range("B16").Formula = "=SMuscleMassMale"

$ws.Range("A17").Value = "Skin"
$ws.Range("B17").Formula = "=SkinMassMale"

$ws.Range("A18").Value = "Organ Mass"
$ws.Range("B18").Formula = "=SUM(B6:B17)"

$ws.Range("A20").Value = "Non-Hepatic Mass"
$ws.Range("B20").Formula = "=SUM(B6:B11,B13:B17)"

$ws.Range("A22").Value = "Hepatic Receptors (/kG BW)"
$ws.Range("B22").Value = 204
$ws.Range("C22").Formula = "=BodyMassMale*B22"
$ws.Range("D22").Value = 12
$ws.Range("E22").Formula = "=0.01*D22*C22"
$ws.Range("F22").Formula = "=C22-E22"

$ws.Range("A23").Value = "Non-Hepatic Receptors (/kG BW)"
$ws.Range("B23").Value = 183
$ws.Range("C23").Formula = "=BodyMassMale*B23"
$ws.Range("D23").Value = 2
$ws.Range("E23").Formula = "=0.01*D23*C23"
$ws.Range("F23").Formula = "=C23-E23"

$ws.Range("C21").Value = "Total (mU)"

$ws.Range("A25").Value = "Secretion=Degradation (mU/Min)"
$ws.Range("B25").Value = 17

$ws.Range("A27").Value = "[Insulin] ECFV"
$ws.Range("B27").Value = 0.02
$ws.Range("C27").Value = 15000
$ws.Range("D27").Formula = "=B27*C27"

$ws.Range("D21").Value = "% Occupied"

$ws.Range("A30").Value = "Degradation"

$ws.Range("A33").Value = "Other"

$ws.Range("B30").Value = "% Total"

$ws.Range("E21").Value = "Occupied"

$ws.Range("F21").Value = "Free"

$ws.Range("C30").Value = "mU/Min"

$ws.Range("D30").Value = "k"

$ws.Range("A28").Value = "[Insulin] Portal Vein"
$ws.Range("B28").Value = 0.052

$ws.Range("A31").Value = "Liver"
$ws.Range("B31").Value = 79
$ws.Range("C31").Formula = "=0.01*B31*B25"
$ws.Range("D31").Formula = "=C31/E22"

$ws.Range("A32").Value = "Kidney"
$ws.Range("B32").Value = 9
$ws.Range("C32").Formula = "=0.01*B32*B25"
$ws.Range("D32").Formula = "=C32/D27"

$ws.Range("A33").Value = "Other"
$ws.Range("B33").Value = 12
$ws.Range("C33").Formula = "=0.01*B33*B25"
$ws.Range("D33").Formula = "=C33/E23"

$ws.Range("C34").Formula = "=SUM(C31:C33)"

# ------------------------------------------------------------------
# 4) Number formatting: integer format for the organ-mass column
#    (matches style index 1 used elsewhere in the workbook).
# ------------------------------------------------------------------
$ws.Range("B6:B18").NumberFormat = "0"
$ws.Range("B20").NumberFormat = "0"

# ------------------------------------------------------------------
# 5) Column widths
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 29.666666666666668
$ws.Columns.Item(3).ColumnWidth = 9.666666666666666
$ws.Columns.Item(4).ColumnWidth = 11.333333333333334
$ws.Columns.Item(5).ColumnWidth = 9.833333333333334

# ------------------------------------------------------------------
# 6) View state: scroll + selection, make this the active sheet/tab
# ------------------------------------------------------------------
$ws.Activate()
$win = $wb.Windows.Item(1)
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("B28").Select()
